$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "goalkeeper pants youth soccer",
    "compression pants men cheap",
    "leggings baseball",
    "indoor knee pads",
    "yoga knee pads 2 pack",
    "compression shorts for basketball",
    "knee pads gym",
    "boys basketball leggings youth",
    "knee sweat pants",
    "knee pads for sports",
    "black softball pants",
    "arthritis hope knee compression sleeve",
    "honeycomb tights",
    "mens sheer pants",
    "mens basketball knee sleeves",
    "boys baseball pants short",
    "spandex work pants men",
    "baseball pants long",
    "capri pants for men",
    "boys basketball tights",
    "tight pants",
    "knee length shorts for men",
    "softball shorts for men",
    "baseball catchers hand pad",
    "cold pad knee",
    "football pants mens with pads",
    "yoga positions chart",
    "hex leg sleeve youth",
    "knees pads yoga",
    "work pants for men with knee pad",
    "tight basketball shorts",
    "adult hockey pants",
    "youth sliding shorts baseball",
    "sheer pants for men",
    "paintball shorts",
    "lacrosse padded shorts",
    "knee pads for paintball",
    "mens big and tall basketball pants",
    "mens shorts below the knee",
    "basketball tights youth",
    "girls sliding shorts softball",
    "youth compression pants boys",
    "boys compression running pants",
    "compression yoga tights",
    "sliding shorts youth girls softball",
    "basketball lot",
    "paintball leg pads",
    "adult black football pants",
    "knee hockey pads",
    "knee sleeve basketball",
    "knee pads work pants",
    "mens knee pads",
    "padded baseball sliding shorts",
    "soccer tights for men",
    "boys tight pants",
    "black leggings for men",
    "men athletic compression pants",
    "youth 5 pad girdle",
    "poc knee pads",
    "men baseball pants black",
    "youth small compression pants",
    "hex knee pads compression leg sleeve",
    "hockey compression pants",
    "softball mens",
    "x compression pants",
    "men compression legging",
    "wrestling knee pads",
    "boys tights for sports youth",
    "hockey pants youth",
    "pant sport men",
    "girls softball sliding pants",
    "boys small compression pants",
    "compression spandex for men",
    "knee pads mtb",
    "sport leggings men",
    "mens sliding shorts",
    "kneeling on the promises",
    "tall baseball pants mens",
    "yoga knee pad thick",
    "men gym leggings",
    "basketball knee sleeves for men",
    "basketball knee sleeves with pads",
    "youth athletic tights",
    "baseball pants adults",
    "tights for soccer",
    "sliding shorts",
    "mens sport leggings",
    "6 inch basketball",
    "football waist pads",
    "mens work pants knee pads",
    "knee sleeve lacrosse",
    "men sliding shorts",
    "padded leggings",
    "yoga pads",
    "mens cold tights",
    "boy tights youth",
    "basketball cycling",
    "basketball knee sleeves youth",
    "soccer compression pants",
    "sliding pads"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
